# PlayerPerformance_3065.xlsx edit:
#  1) "ODI Batting Extra": the empty placeholder cells C2:E2 and C3:E3 are
#     removed (cleared to blank) - they used to hold empty inline strings.
#  2) A new worksheet "ODI Bowling Extra" (sheetId 5) is appended right
#     after "ODI Batting Extra", carrying MATCH_CODE / MAIDEN_OVERS /
#     PERCENT_WICKETS_OF_ALL columns for the bowling-extra scrape.

$wb = $excel.ActiveWorkbook

# --- 1. Trim the stray empty placeholder cells on "ODI Batting Extra" ---
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$battingExtra.Range("C2:E3").ClearContents()

# --- 2. Add the new "ODI Bowling Extra" sheet right after it ---
$bowlingExtra = $wb.Worksheets.Add($null, $battingExtra)
$bowlingExtra.Name = "ODI Bowling Extra"

$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")

$rows = @(
  @("3637", "0", "20.00%"),
  @("3638", "1", "10.00%"),
  @("3639", $null, $null),
  @("3640", $null, $null),
  @("3641", $null, $null),
  @("3642", $null, $null),
  @("3664", "1", $null),
  @("3667", $null, $null),
  @("3670", "1", $null),
  @("3672", "0", "20.00%"),
  @("3735", "2", "40.00%"),
  @("3738", $null, $null),
  @("3744", "1", "10.00%"),
  @("3746", $null, $null),
  @("3749", "0", $null),
  @("3756", "0", $null),
  @("3761", $null, $null),
  @("3769", "0", $null),
  @("3780", "1", "20.00%"),
  @("3785", "0", "10.00%")
)

# Force every cell in the used range to Text so numeric-looking values
# (match codes, "0"/"1" maiden-over counts, "20.00%" strings) are stored
# as literal text, matching the scraper's inlineStr output.
$usedRange = $bowlingExtra.Range("A1:C21")
$usedRange.NumberFormat = "@"

for ($col = 1; $col -le 3; $col++) {
    $bowlingExtra.Cells.Item(1, $col).Value = $headers[$col - 1]
}

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($col = 1; $col -le 3; $col++) {
        $val = $rowData[$col - 1]
        if ($val -ne $null) {
            $bowlingExtra.Cells.Item($r, $col).Value = $val
        }
    }
}

# Match the bold / bordered / centered-top header styling used by every
# other sheet in the workbook.
$hdr = $bowlingExtra.Range("A1:C1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

$bowlingExtra.Range("A1").Select()
